# Applies the "Copy Input files from P drive. Delete old input files." edit
# to the Constants sheet of the STS IR Bot Performer Config workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")

# --- Insert a new row above the "PathOtherDeductionsList" row (row 17) ---
# This shifts the old rows 17..116 down to 18..117 and grows the used range
# by one row (766 -> 767), matching the target dimension change.
$ws.Rows.Item(17).Insert()

# --- New row 17: ListInputFiles (brand new constant), highlighted green ---
$ws.Range("A17").Value = "ListInputFiles"
$ws.Range("B17").Value = "PathOtherDeductionsList,PathMarginsList,PathCustomerNameList"
$ws.Range("C17").Value = "Constant names (from this sheet) of each input file that must be copied from the P drive"
$ws.Range("A17:C17").Interior.Color = 5296274   # RGB(80,208,146) == BGR 0x50D092 -> fgColor FF92D050

# --- Row 18 (was old row 17): PathOtherDeductionsList now points at the new
# local destination path ("Data\Other Deductions List.xlsx" instead of
# "Data\Input\Other Deductions List.xlsx"). Highlighted yellow. ---
$ws.Range("B18").Value = "Data\Other Deductions List.xlsx"
$ws.Range("A18:B18").Interior.Color = 65535     # RGB(255,255,0) yellow

# --- Row 19 (was old row 18): PathPDriveFolder now points at the new
# "InputFiles" share; the previous TaxSolver Files path is preserved in
# column C for reference. Highlighted green. ---
$ws.Range("C19").Value = $ws.Range("B19").Value()
$ws.Range("B19").Value = "\\somproddfs1.prod.sovos.org\depts\TaxReturnOutSourcing\Preparer\UIPathPublish\IR Bot Temp Files\InputFiles"
$ws.Range("A19:C19").Interior.Color = 5296274

# --- Row 20 (was old row 19): PathTemplatesFile - unchanged, no highlight ---

# --- Row 21 (was old row 20): PathMarginsList - now part of the copied
# input-file list, highlighted yellow, values unchanged. ---
$ws.Range("A21:B21").Interior.Color = 65535

# --- Row 22 (was old row 21): PathCustomerNameList - same treatment. ---
$ws.Range("A22:B22").Interior.Color = 65535

# --- Update the view / selection to match the saved state ---
[void]$ws.Activate()
$ws.Range("B21").Select()
